$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("header")
$ws.Range("A1").Value = "orcid: https://orcid.org/"
$ws.Range("A2").Value = "sssom: https://w3id.org/sssom/"
$ws.Range("A3").Value = "semapv: https://w3id.org/semapv/vocab/"
$ws.Range("A4").Value = "crosswalk: https://w3id.org/env/neap/crosswalk/"
$ws.Range("A5").Value = "status: https://w3id.org/env/neap/status/"
$ws.Range("A6").Value = "get: https://global-ecosystems.org/explore/"
$ws.Range("A7").Value = "smartline: https://w3id.org/env/neap/smartline/"
$ws.Range("A8").Value = "map: http://w3id.org/env/neap/smartline-get/"
$ws.Range("A9").ClearContents()
$ws.Range("A10").ClearContents()
